$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column D: "primary_key" header plus "<Unspecified>" values for
# each of the 12 data rows (rows 2-13).
$ws.Range("D1").Value = "primary_key"

for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 4).Value = "<Unspecified>"
}
